$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Formula Samples")

# New row 12: label, computed formula (with extra algebraic calc), explanation
$ws.Range("A12").Value = "Sum with extra algebraic calculation"
$ws.Range("B12").Formula = "=B2+B3+B4+4-2"
$ws.Range("C12").Value = "> This cell will also be considered as right, since the extra algebraic calculation doesn't happen inside custom function - allowing it to be simplified."

# Column A needs to widen to fit the new, longer label (best-fit width)
$ws.Columns.Item(1).ColumnWidth = 30.65

# Match the saved selection state from the diff (active cell now B7)
$ws.Range("B7").Select()
